# Daily attendance processing - 2025-10-12 03:36:40
# Normalize the "Recorded By" column (G) so that the automated "System"
# actor is always listed first, ahead of the human/backup account that
# triggered the recording (e.g. "user@x.com, System" -> "System, user@x.com").
# The "admin@admin.com" account is left untouched / already conforms.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $text = $cell.Value2

    if ($null -eq $text) { continue }

    $parts = $text -split ', '
    if ($parts.Count -lt 2) { continue }

    $first = $parts[0].Trim()
    $second = $parts[1].Trim()

    if ($second -eq 'System' -and $first -ne 'System' -and ($first -eq 'dnasr281@gmail.com' -or $first -eq 'backup@backdoor.com')) {
        $newParts = @('System', $first)
        if ($parts.Count -gt 2) {
            $rest = $parts[2..($parts.Count - 1)]
            $newParts = $newParts + $rest
        }
        $cell.Value = [string]::Join(', ', $newParts)
    }
}
